# Insert two new rows at 1125, pushing existing rows 1125:1219 down to 1127:1221,
# then populate the two newly-inserted rows with their values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1125:1126").Insert()

# New row 1125
$ws.Range("A1125").Value2 = 10
$ws.Range("B1125").Value2 = "Vega Modelo de Temuco"
$ws.Range("C1125").Value2 = "La Araucanía"
$ws.Range("D1125").Value2 = 45021
$ws.Range("E1125").Value2 = 9
$ws.Range("F1125").Value2 = 100112021
$ws.Range("G1125").Value2 = "Ají"
$ws.Range("H1125").Value2 = "Americana (o)"
$ws.Range("I1125").Value2 = "Primera"
$ws.Range("J1125").Value2 = 35
$ws.Range("K1125").Value2 = 20000
$ws.Range("L1125").Value2 = 20000
$ws.Range("M1125").Value2 = 20000
$ws.Range("N1125").Value2 = '$/caja 25 kilos'
$ws.Range("O1125").Value2 = "Región del Maule"
$ws.Range("P1125").Value2 = 800
$ws.Range("Q1125").Value2 = 25
$ws.Range("R1125").Value2 = "Hortaliza"

# New row 1126
$ws.Range("A1126").Value2 = 10
$ws.Range("B1126").Value2 = "Vega Modelo de Temuco"
$ws.Range("C1126").Value2 = "La Araucanía"
$ws.Range("D1126").Value2 = 45021
$ws.Range("E1126").Value2 = 9
$ws.Range("F1126").Value2 = 100112021
$ws.Range("G1126").Value2 = "Ají"
$ws.Range("H1126").Value2 = "Americana (o)"
$ws.Range("I1126").Value2 = "Primera"
$ws.Range("J1126").Value2 = 55
$ws.Range("K1126").Value2 = 20000
$ws.Range("L1126").Value2 = 20000
$ws.Range("M1126").Value2 = 20000
$ws.Range("N1126").Value2 = '$/saco 25 kilos'
$ws.Range("O1126").Value2 = "Región del Maule"
$ws.Range("P1126").Value2 = 800
$ws.Range("Q1126").Value2 = 25
$ws.Range("R1126").Value2 = "Hortaliza"
